# Auto-generated Excel COM-interop script applying the scraped diff.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across
# several rows on the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets, matching
# the "chore: update Sheets via scheduled runner" commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 349.5
$ws.Range("J45").Value = 449
$ws.Range("L45").Value = 1347
$ws.Range("N45").Value = -1731

$ws.Range("H98").Value = 2181.973
$ws.Range("I98").Value = 2242.4167
$ws.Range("J98").Value = 6
$ws.Range("K98").Value = 2242.4167
$ws.Range("L98").Value = 6
$ws.Range("M98").Value = -744.4167000000002
$ws.Range("N98").Value = -3002

$ws.Range("H106").Value = 6165.3335
$ws.Range("I106").Value = 6085.5654
$ws.Range("K106").Value = 6085.5654
$ws.Range("M106").Value = -5454.5654

$ws.Range("H115").Value = 285.45456
$ws.Range("I115").Value = 285.45456
$ws.Range("K115").Value = 856.36368
$ws.Range("M115").Value = 710.63632

$ws.Range("H122").Value = 2181.973
$ws.Range("I122").Value = 2242.4167
$ws.Range("J122").Value = 6
$ws.Range("K122").Value = 6727.250100000001
$ws.Range("L122").Value = 18
$ws.Range("M122").Value = -4277.250100000001
$ws.Range("N122").Value = -4918


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 148962.14
$ws.Range("I45").Value = 289641.44
$ws.Range("J45").Value = 8282.857
$ws.Range("K45").Value = 289641.44
$ws.Range("L45").Value = 8282.857
$ws.Range("M45").Value = -289264.44
$ws.Range("N45").Value = -9036.857

$ws.Range("H61").Value = 8836.087
$ws.Range("I61").Value = 7622.533
$ws.Range("K61").Value = 7622.533
$ws.Range("M61").Value = -7410.533

$ws.Range("H136").Value = 8836.087
$ws.Range("I136").Value = 7622.533
$ws.Range("K136").Value = 22867.599
$ws.Range("M136").Value = -20317.599


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 70000
$ws.Range("I53").Value = 70000
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 70000
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -69426
$ws.Range("N53").Value = $null

$ws.Range("H99").Value = 2853.2195
$ws.Range("I99").Value = 1587.3928
$ws.Range("K99").Value = 1587.3928
$ws.Range("M99").Value = -89.39280000000008

$ws.Range("H105").Value = 3381.3333
$ws.Range("I105").Value = 2900.5881
$ws.Range("J105").Value = 5424.5
$ws.Range("K105").Value = 2900.5881
$ws.Range("L105").Value = 5424.5
$ws.Range("M105").Value = -1153.5881
$ws.Range("N105").Value = -8918.5

$ws.Range("H134").Value = 7472.6113
$ws.Range("I134").Value = 7469.3125
$ws.Range("J134").Value = 7499
$ws.Range("K134").Value = 22407.9375
$ws.Range("L134").Value = 22497
$ws.Range("M134").Value = -19872.9375
$ws.Range("N134").Value = -27567


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3757.4
$ws.Range("I16").Value = 2650.3635
$ws.Range("K16").Value = 2650.3635
$ws.Range("M16").Value = -2363.3635

$ws.Range("H31").Value = 4013.2144
$ws.Range("I31").Value = 3395.6
$ws.Range("J31").Value = 4574.6816
$ws.Range("K31").Value = 3395.6
$ws.Range("L31").Value = 4574.6816
$ws.Range("M31").Value = -3100.6
$ws.Range("N31").Value = -5164.6816

$ws.Range("H34").Value = 4013.2144
$ws.Range("I34").Value = 3395.6
$ws.Range("J34").Value = 4574.6816
$ws.Range("K34").Value = 3395.6
$ws.Range("L34").Value = 4574.6816
$ws.Range("M34").Value = -3193.6
$ws.Range("N34").Value = -4978.6816

$ws.Range("H58").Value = 5699.1763
$ws.Range("I58").Value = 3129.6
$ws.Range("K58").Value = 3129.6
$ws.Range("M58").Value = -2926.6

$ws.Range("H94").Value = 1743.6364
$ws.Range("I94").Value = 966
$ws.Range("J94").Value = 2035.25
$ws.Range("K94").Value = 966
$ws.Range("L94").Value = 2035.25
$ws.Range("M94").Value = -515
$ws.Range("N94").Value = -2937.25

$ws.Range("H113").Value = 3757.4
$ws.Range("I113").Value = 2650.3635
$ws.Range("K113").Value = 2650.3635
$ws.Range("M113").Value = -480.3634999999999

$ws.Range("H122").Value = 4653.1333
$ws.Range("I122").Value = 4449.7856
$ws.Range("K122").Value = 13349.3568
$ws.Range("M122").Value = -10899.3568

$ws.Range("H132").Value = 4595
$ws.Range("I132").Value = 2868.3333
$ws.Range("K132").Value = 8604.999899999999
$ws.Range("M132").Value = -6074.999899999999

$ws.Range("I134").Value = 2756.2666
$ws.Range("J134").Value = 10352.429
$ws.Range("K134").Value = 8268.799800000001
$ws.Range("L134").Value = 31057.287
$ws.Range("M134").Value = -5733.799800000001
$ws.Range("N134").Value = -36127.287

$ws.Range("H136").Value = 5699.1763
$ws.Range("I136").Value = 3129.6
$ws.Range("K136").Value = 9388.799999999999
$ws.Range("M136").Value = -6838.799999999999

$ws.Range("H140").Value = 99017.8
$ws.Range("J140").Value = 99017.8
$ws.Range("L140").Value = 99017.8
$ws.Range("N140").Value = -109377.8


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 865.8
$ws.Range("I8").Value = 865.8
$ws.Range("K8").Value = 2597.4
$ws.Range("M8").Value = -2458.4

$ws.Range("H12").Value = 22.125
$ws.Range("I12").Value = 23.666666
$ws.Range("K12").Value = 70.99999800000001
$ws.Range("M12").Value = 102.000002

$ws.Range("H56").Value = 22201.363
$ws.Range("I56").Value = 22201.363
$ws.Range("K56").Value = 22201.363
$ws.Range("M56").Value = -21671.363

$ws.Range("H92").Value = 197.5
$ws.Range("J92").Value = 197.5
$ws.Range("L92").Value = 592.5
$ws.Range("N92").Value = -3088.5

$ws.Range("H109").Value = 1757.4286
$ws.Range("I109").Value = 745.3333
$ws.Range("J109").Value = 2516.5
$ws.Range("K109").Value = 2235.9999
$ws.Range("L109").Value = 7549.5
$ws.Range("M109").Value = -1195.9999
$ws.Range("N109").Value = -9629.5

$ws.Range("H131").Value = 31251976
$ws.Range("I131").Value = 71429160
$ws.Range("J131").Value = 3055.3333
$ws.Range("K131").Value = 214287480
$ws.Range("L131").Value = 9165.999899999999
$ws.Range("M131").Value = -214282440
$ws.Range("N131").Value = -19245.9999


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2589
$ws.Range("I102").Value = 1963.0588
$ws.Range("J102").Value = 5249.25
$ws.Range("K102").Value = 1963.0588
$ws.Range("L102").Value = 5249.25
$ws.Range("M102").Value = -341.0588
$ws.Range("N102").Value = -8493.25

$ws.Range("H126").Value = 7885.3335
$ws.Range("I126").Value = 7885.3335
$ws.Range("K126").Value = 23656.0005
$ws.Range("M126").Value = -21186.0005


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 54333.332
$ws.Range("I74").Value = 54333.332
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 54333.332
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -53335.332
$ws.Range("N74").Value = $null

$ws.Range("H77").Value = 54333.332
$ws.Range("I77").Value = 54333.332
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 162999.996
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -158007.996
$ws.Range("N77").Value = $null

$ws.Range("H132").Value = 18567.715
$ws.Range("I132").Value = 18567.715
$ws.Range("K132").Value = 55703.145
$ws.Range("M132").Value = -53173.145

$ws.Range("H136").Value = 3451.8823
$ws.Range("I136").Value = 3054
$ws.Range("K136").Value = 9162
$ws.Range("M136").Value = -6612


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 31296.666
$ws.Range("J54").Value = 31296.666
$ws.Range("L54").Value = 31296.666
$ws.Range("N54").Value = -32336.666

$ws.Range("H126").Value = 6085.2354
$ws.Range("I126").Value = 3787.1667
$ws.Range("K126").Value = 11361.5001
$ws.Range("M126").Value = -8891.500100000001

$ws.Range("H136").Value = 6248.36
$ws.Range("I136").Value = 5409.143
$ws.Range("K136").Value = 16227.429
$ws.Range("M136").Value = -13677.429
